$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "...will be a web-based accessibility-focused..." ->
#    "...will be a in-house application accessibility-focused..."
#    (and split into 3 runs, matching the target XML exactly)
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("web-based", $false, $false, $false, $false, $false, $true, 1, $false, "in-house application", 2)
if ($found) {
    # Force the newly-inserted text to live in its own run (distinct rPr
    # object) instead of being silently re-merged with its neighbours.
    $rng.Font.Size = 13
    $rng.Font.Size = 12
}

# ------------------------------------------------------------------
# 2) Merge the "...Heroku, PythonAnywhere, or [Vercel] for frontend"
#    run/proofErr/run/proofErr/run sequence into a single run.
# ------------------------------------------------------------------
$rng2 = $d.Content
$hostingText = "- Hosting: MMU internal servers, or cloud platforms like Heroku, PythonAnywhere, or Vercel for frontend"
$rng2.Find.Execute($hostingText, $false, $false, $false, $false, $false, $true, 1, $false, $hostingText, 2) | Out-Null

# ------------------------------------------------------------------
# 3) Merge the "- [AccessMap] (University of Washington): ..." run/
#    proofErr/run/proofErr/run sequence into a single run.
# ------------------------------------------------------------------
$rng3 = $d.Content
$accessMapText = "- AccessMap (University of Washington): Offers route planning focused on sidewalk slopes and mobility in urban settings."
$rng3.Find.Execute($accessMapText, $false, $false, $false, $false, $false, $true, 1, $false, $accessMapText, 2) | Out-Null

# ------------------------------------------------------------------
# 4) Subtitle style: add a first-line indent of 720 twips (36pt).
# ------------------------------------------------------------------
$subtitle = $d.Styles("Subtitle")
$subtitle.ParagraphFormat.FirstLineIndent = 36
